$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2935490919770132
$ws.Range("C2").Value = 0.04457165816555175
$ws.Range("D2").Value = 0.03334822985983266
$ws.Range("F2").Value = 0.5500782417717929
$ws.Range("G2").Value = 0.002413447129026704
$ws.Range("I2").Value = 0.4222757950674723
$ws.Range("K2").Value = 0.2990256334740593
$ws.Range("N2").Value = 1.189236620751288
$ws.Range("O2").Value = 1.857100913325638

$ws.Range("B3").Value = 0.2585797311558338
$ws.Range("C3").Value = 0.03959081992960023
$ws.Range("D3").Value = 0.03002597117231431
$ws.Range("F3").Value = 0.5519381186377856
$ws.Range("G3").Value = 0.002415292165664462
$ws.Range("I3").Value = 0.4274429606471148
$ws.Range("K3").Value = 0.2615162650164677
$ws.Range("N3").Value = 1.199097099743849
$ws.Range("O3").Value = 1.873458363975587

$ws.Range("B4").Value = 0.2370770665903308
$ws.Range("C4").Value = 0.03651435465314989
$ws.Range("D4").Value = 0.02797201785985237
$ws.Range("F4").Value = 0.5534368788419073
$ws.Range("G4").Value = 0.002416486112989779
$ws.Range("I4").Value = 0.430852329168431
$ws.Range("K4").Value = 0.2384247619132083
$ws.Range("N4").Value = 1.205565533908441
$ws.Range("O4").Value = 1.884619283030048

$ws.Range("B5").Value = 0.2283073025316469
$ws.Range("C5").Value = 0.03525616233703488
$ws.Range("D5").Value = 0.0271315237780172
$ws.Range("F5").Value = 0.5541373771130296
$ws.Range("G5").Value = 0.002416988060423152
$ws.Range("I5").Value = 0.4323011847498197
$ws.Range("K5").Value = 0.2290001859672941
$ws.Range("N5").Value = 1.208305678571108
$ws.Range("O5").Value = 1.889448394448493

$ws.Range("B6").Value = 0.2268506702042998
$ws.Range("C6").Value = 0.03504697030781756
$ws.Range("D6").Value = 0.02699175089325223
$ws.Range("F6").Value = 0.5542591151725276
$ws.Range("G6").Value = 0.002417072340332862
$ws.Range("I6").Value = 0.4325453600367783
$ws.Range("K6").Value = 0.2274343791988969
$ws.Range("N6").Value = 1.208766973819781
$ws.Range("O6").Value = 1.890267232375194

$ws.Range("B7").Value = 0.2369588230255602
$ws.Range("C7").Value = 0.03649740438989113
$ws.Range("D7").Value = 0.02796069673581059
$ws.Range("F7").Value = 0.5534459626152923
$ws.Range("G7").Value = 0.002416492820308932
$ws.Range("I7").Value = 0.4308716280022971
$ws.Range("K7").Value = 0.2382977170606324
$ws.Range("N7").Value = 1.205602066418571
$ws.Range("O7").Value = 1.884683272610332

$ws.Range("B8").Value = 0.281498525036028
$ws.Range("C8").Value = 0.04285808615233577
$ws.Range("D8").Value = 0.03220566343152598
$ws.Range("F8").Value = 0.5506454892541015
$ws.Range("G8").Value = 0.002414070644902575
$ws.Range("I8").Value = 0.4240082982276583
$ws.Range("K8").Value = 0.2861053517272865
$ws.Range("N8").Value = 1.192550627687559
$ws.Range("O8").Value = 1.862509033383077

$ws.Range("B9").Value = 0.3685687087306064
$ws.Range("C9").Value = 0.05518428851182478
$ws.Range("D9").Value = 0.04041664521025723
$ws.Range("F9").Value = 0.547984267417938
$ws.Range("G9").Value = 0.002409803484515777
$ws.Range("I9").Value = 0.4124279624439708
$ws.Range("K9").Value = 0.3793514906310804
$ws.Range("N9").Value = 1.170238259266817
$ws.Range("O9").Value = 1.827893189477791

$ws.Range("B10").Value = 0.4323479115265627
$ws.Range("C10").Value = 0.06414803387042411
$ws.Range("D10").Value = 0.0463783453306803
$ws.Range("F10").Value = 0.5477548853551895
$ws.Range("G10").Value = 0.002406959897230769
$ws.Range("I10").Value = 0.4050656851707721
$ws.Range("K10").Value = 0.4475264653480622
$ws.Range("N10").Value = 1.155840378973146
$ws.Range("O10").Value = 1.807870475696276

$ws.Range("B11").Value = 0.4613161491629967
$ws.Range("C11").Value = 0.06820531245364236
$ws.Range("D11").Value = 0.04907472939170532
$ws.Range("F11").Value = 0.5480253861708348
$ws.Range("G11").Value = 0.002405728980213843
$ws.Range("I11").Value = 0.4019653067373596
$ws.Range("K11").Value = 0.4784637044757574
$ws.Range("N11").Value = 1.149722489756314
$ws.Range("O11").Value = 1.799937080743561

$ws.Range("B12").Value = 0.4722785977968726
$ws.Range("C12").Value = 0.06973870583503583
$ws.Range("D12").Value = 0.05009349086955694
$ws.Range("F12").Value = 0.5481817118637764
$ws.Range("G12").Value = 0.00240527182791948
$ws.Range("I12").Value = 0.4008270606173916
$ws.Range("K12").Value = 0.4901673469138359
$ws.Range("N12").Value = 1.147467818691993
$ws.Range("O12").Value = 1.797101940213224

$ws.Range("B13").Value = 0.4699179692720463
$ws.Range("C13").Value = 0.06940859728059934
$ws.Range("D13").Value = 0.04987418554155454
$ws.Range("F13").Value = 0.5481456475958453
$ws.Range("G13").Value = 0.002405369885451027
$ws.Range("I13").Value = 0.4010706096198611
$ws.Range("K13").Value = 0.4876472860681247
$ws.Range("N13").Value = 1.147950643915635
$ws.Range("O13").Value = 1.797705017799601

$ws.Range("B14").Value = 0.4622181845377042
$ws.Range("C14").Value = 0.06833152643385176
$ws.Range("D14").Value = 0.04915858993018674
$ws.Range("F14").Value = 0.5480371671456936
$ws.Range("G14").Value = 0.002405691190600271
$ws.Range("I14").Value = 0.4018709448967641
$ws.Range("K14").Value = 0.4794268070725423
$ws.Range("N14").Value = 1.149535753491001
$ws.Range("O14").Value = 1.799700443195547

$ws.Range("B15").Value = 0.4575008867189467
$ws.Range("C15").Value = 0.06767139507496722
$ws.Range("D15").Value = 0.04871996573632487
$ws.Range("F15").Value = 0.547977737874632
$ws.Range("G15").Value = 0.002405889166104409
$ws.Range("I15").Value = 0.4023658364166636
$ws.Range("K15").Value = 0.4743899923327035
$ws.Range("N15").Value = 1.150514757016495
$ws.Range("O15").Value = 1.800944718768847

$ws.Range("B16").Value = 0.4304537998334013
$ws.Range("C16").Value = 0.0638824643393292
$ws.Range("D16").Value = 0.04620181129507728
$ws.Range("F16").Value = 0.547744749580815
$ws.Range("G16").Value = 0.002407041597940672
$ws.Range("I16").Value = 0.4052733103112338
$ws.Range("K16").Value = 0.4455030565208347
$ws.Range("N16").Value = 1.156248881725048
$ws.Range("O16").Value = 1.808412594643073

$ws.Range("B17").Value = 0.413849210709202
$ws.Range("C17").Value = 0.06155280270046148
$ws.Range("D17").Value = 0.04465296521908613
$ws.Range("F17").Value = 0.5476978197700006
$ws.Range("G17").Value = 0.002407764594585449
$ws.Range("I17").Value = 0.4071206849023952
$ws.Range("K17").Value = 0.4277619163696897
$ws.Range("N17").Value = 1.159877135269937
$ws.Range("O17").Value = 1.813294918201422

$ws.Range("B18").Value = 0.4042944730280169
$ws.Range("C18").Value = 0.06021093002752309
$ws.Range("D18").Value = 0.04376064262885393
$ws.Range("F18").Value = 0.5477061100987584
$ws.Range("G18").Value = 0.002408186341605349
$ws.Range("I18").Value = 0.4082066585233761
$ws.Range("K18").Value = 0.417550579950273
$ws.Range("N18").Value = 1.162004654568456
$ws.Range("O18").Value = 1.81621369616704

$ws.Range("B19").Value = 0.401058699696506
$ws.Range("C19").Value = 0.05975626889554064
$ws.Range("D19").Value = 0.04345826716863144
$ws.Range("F19").Value = 0.5477149772487024
$ws.Range("G19").Value = 0.002408330152320549
$ws.Range("I19").Value = 0.4085783707917052
$ws.Range("K19").Value = 0.4140919986084839
$ws.Range("N19").Value = 1.162731978173653
$ws.Range("O19").Value = 1.817220936037401

$ws.Range("B20").Value = 0.4156172390185873
$ws.Range("C20").Value = 0.06180099771914627
$ws.Range("D20").Value = 0.04481799477794368
$ws.Range("F20").Value = 0.5476991639258486
$ws.Range("G20").Value = 0.002407687020034088
$ws.Range("I20").Value = 0.4069216053965654
$ws.Range("K20").Value = 0.4296512311408094
$ws.Range("N20").Value = 1.159486695360471
$ws.Range("O20").Value = 1.812763739705403

$ws.Range("B21").Value = 0.4644799986984367
$ws.Range("C21").Value = 0.06864797057968985
$ws.Range("D21").Value = 0.04936884057775615
$ws.Range("F21").Value = 0.5480675679293654
$ws.Range("G21").Value = 0.00240559657230638
$ws.Range("I21").Value = 0.4016348951985087
$ws.Range("K21").Value = 0.4818416819609297
$ws.Range("N21").Value = 1.149068484991481
$ws.Range("O21").Value = 1.799109749293308

$ws.Range("B22").Value = 0.4963723835095095
$ws.Range("C22").Value = 0.07310528155139195
$ws.Range("D22").Value = 0.05232964626880232
$ws.Range("F22").Value = 0.5486224649527003
$ws.Range("G22").Value = 0.002404282604999409
$ws.Range("I22").Value = 0.3983884068823755
$ws.Range("K22").Value = 0.5158831095178016
$ws.Range("N22").Value = 1.14262118251618
$ws.Range("O22").Value = 1.79117155399797

$ws.Range("B23").Value = 0.4793549112376922
$ws.Range("C23").Value = 0.07072796601600828
$ws.Range("D23").Value = 0.05075065575985604
$ws.Range("F23").Value = 0.548297567163047
$ws.Range("G23").Value = 0.002404979125061675
$ws.Range("I23").Value = 0.4001020149059933
$ws.Range("K23").Value = 0.497721018264059
$ws.Range("N23").Value = 1.146029158217296
$ws.Range("O23").Value = 1.795318116196484

$ws.Range("B24").Value = 0.4148179399627452
$ws.Range("C24").Value = 0.06168879659279014
$ws.Range("D24").Value = 0.04474339073628641
$ws.Range("F24").Value = 0.5476984463673702
$ws.Range("G24").Value = 0.002407722072461139
$ws.Range("I24").Value = 0.4070115348179435
$ws.Range("K24").Value = 0.4287971083813318
$ws.Range("N24").Value = 1.159663083696522
$ws.Range("O24").Value = 1.813003537061732

$ws.Range("B25").Value = 0.3450456683714833
$ws.Range("C25").Value = 0.05186574397218635
$ws.Range("D25").Value = 0.03820767978419326
$ws.Range("F25").Value = 0.5484011425458775
$ws.Range("G25").Value = 0.002410906475056687
$ws.Range("I25").Value = 0.4153595742022791
$ws.Range("K25").Value = 0.3541825752969885
$ws.Range("N25").Value = 1.175923614966436
$ws.Range("O25").Value = 1.836307967720614

